$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the stack-trace line "...finish(InstanceParser.java:177)"
# was split across two runs by a "_GoBack" bookmark; merge it back into
# a single run and drop the bookmark from this spot.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("at com.novaordis.em.ec2.parser.InstanceParser.finish(Ins", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $r1.Start

$r2 = $d.Content
$r2.Find.Execute("tanceParser.java:177)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeEnd = $r2.End

$mergeLen = $mergeEnd - $mergeStart
$finalLine = "at com.novaordis.em.ec2.parser.InstanceParser.finish(InstanceParser.java:177)"

# Route the merge through a placeholder string first: the target text is
# byte-identical to the two runs concatenated, and a same-text assignment
# is treated as a no-op (bookmark would survive). Changing it, then
# changing it back, forces the runs/bookmark to actually be rewritten.
$placeholder = $d.Range($mergeStart, $mergeEnd)
$placeholder.Text = "PLACEHOLDER_MERGE_TOKEN"
$final = $d.Range($mergeStart, $mergeStart + ("PLACEHOLDER_MERGE_TOKEN".Length))
$final.Text = $finalLine

# ---------------------------------------------------------------------
# Change 2: drop the "em version does not report the correct version"
# bullet paragraph entirely (together with the blank paragraph that used
# to trail it becoming the sole survivor), and plant the "_GoBack"
# bookmark, on its own, in that surviving empty paragraph.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like '*em version*does not report the correct version*') {
        $target = $i
        break
    }
}

$bulletPara = $d.Paragraphs($target)
$bulletRange = $d.Range($bulletPara.Range.Start, $bulletPara.Range.End)
$bulletRange.Delete()

# The paragraph that used to follow the bullet (an empty paragraph with
# no runs at all) is now at the same index. Bookmarking directly into a
# run-less paragraph silently lands the bookmark at document position 0,
# so give it transient text first, bookmark that, then remove the text.
$emptyPara = $d.Paragraphs($target)
$emptyStart = $emptyPara.Range.Start
$stub = $d.Range($emptyStart, $emptyStart)
$stub.InsertBefore("TEMP_STUB")

$bmAnchor = $d.Paragraphs($target).Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmAnchor, $bmAnchor)) | Out-Null

$stubRange = $d.Range($bmAnchor, $bmAnchor + ("TEMP_STUB".Length))
$stubRange.Text = ""
